$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$src = $ws.Range("D1")
$src.Copy()
$ws.Range("A9").PasteSpecial(-4122)
